$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows above the current row 2 (pushing all existing
# data rows down by 6), then strip the formatting that Insert() copies
# down from the header row so the new rows stay unstyled like the other
# data rows.
$ws.Rows("2:7").Insert()
$ws.Range("A2:D7").Style = "Normal"

# Populate the 6 newly inserted data rows.
$ws.Cells.Item(2, 1).Value = "120 Racecourse Rd, Flemington VIC 3031"
$ws.Cells.Item(2, 2).Value = -37.788414
$ws.Cells.Item(2, 3).Value = 144.936951
$ws.Cells.Item(2, 4).Value = "Moonee Valley (C)"

$ws.Cells.Item(3, 1).Value = "30A The Centreway, Preston VIC 3072"
$ws.Cells.Item(3, 2).Value = -37.73908
$ws.Cells.Item(3, 3).Value = 145.002236
$ws.Cells.Item(3, 4).Value = "Darebin (C)"

$ws.Cells.Item(4, 1).Value = "100 Victoria St, Flemington VIC 3031"
$ws.Cells.Item(4, 2).Value = -37.785949
$ws.Cells.Item(4, 3).Value = 144.935308
$ws.Cells.Item(4, 4).Value = "Moonee Valley (C)"

$ws.Cells.Item(5, 1).Value = "Arrival Dr, Melbourne Airport VIC 3045"
$ws.Cells.Item(5, 2).Value = -37.670681
$ws.Cells.Item(5, 3).Value = 144.850306
$ws.Cells.Item(5, 4).Value = "Hume (C)"

$ws.Cells.Item(6, 1).Value = "373 Collins St, Melbourne VIC 3000"
$ws.Cells.Item(6, 2).Value = -37.817011
$ws.Cells.Item(6, 3).Value = 144.962265
$ws.Cells.Item(6, 4).Value = "Melbourne (C)"

$ws.Cells.Item(7, 1).Value = "33 Alfred St, North Melbourne VIC 3051"
$ws.Cells.Item(7, 2).Value = -37.791989
$ws.Cells.Item(7, 3).Value = 144.940508
$ws.Cells.Item(7, 4).Value = "Melbourne (C)"
